$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EPBDS-11426: Not all java syntax formats are supported for Parameter types
# in Rules tables. Update the sample parameter declaration from the
# unsupported "Integer [] intArr" syntax to the supported "Integer [1] intArr"
# syntax.
$ws.Range("C6").Value = "Integer [1] intArr"

# Move/leave the selection where Excel left it after the edit.
$ws.Range("E18").Select()
